# Update cryptocurrency price/volume data per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.160.42'
$ws.Range('E2').Value = '  -5.24%  '
$ws.Range('D3').Value = '3.380.48'
$ws.Range('E3').Value = '  -7.02%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = "'184.97"
$ws.Range('E5').Value = '  -9.26%  '
$ws.Range('D6').Value = "'527.95"
$ws.Range('E6').Value = '  -8.17%  '
$ws.Range('D7').Value = "'0.599"
$ws.Range('E7').Value = '  -4.34%  '
$ws.Range('D8').Value = '3.376.49'
$ws.Range('E8').Value = '  -6.95%  '
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('D10').Value = "'0.624"
$ws.Range('E10').Value = '  -9.69%  '
$ws.Range('D11').Value = "'57.43"
$ws.Range('E11').Value = '  -7.01%  '
$ws.Range('E12').Value = '  -13.28%  '
$ws.Range('D13').Value = "'0.0000253"
$ws.Range('E13').Value = '  -12.61%  '
$ws.Range('D14').Value = "'9.25"
$ws.Range('E14').Value = '  -9.25%  '
$ws.Range('D15').Value = '3.904.26'
$ws.Range('E15').Value = '  -7.26%  '
$ws.Range('E16').Value = '  -3.57%  '
$ws.Range('D17').Value = '3.371.16'
$ws.Range('E17').Value = '  -7.33%  '
$ws.Range('D18').Value = '64.781.15'
$ws.Range('E18').Value = '  -5.55%  '
$ws.Range('D19').Value = "'17.37"
$ws.Range('E19').Value = '  -9.03%  '
$ws.Range('D20').Value = "'11.04"
$ws.Range('E20').Value = '  -11.93%  '
$ws.Range('D21').Value = "'0.965"
$ws.Range('E21').Value = '  -11.05%  '
$ws.Range('D22').Value = "'371.82"
$ws.Range('E22').Value = '  -9.07%  '
$ws.Range('D23').Value = "'80.93"
$ws.Range('E23').Value = '  -6.21%  '
$ws.Range('D24').Value = "'3.72"
$ws.Range('E24').Value = '  -12.32%  '
$ws.Range('D25').Value = "'10.84"
$ws.Range('E25').Value = '  -16.44%  '
$ws.Range('D26').Value = "'3.76"
$ws.Range('E26').Value = '  -5.73%  '
$ws.Range('E27').Value = '  -5.44%  '
$ws.Range('D28').Value = "'2.64"
$ws.Range('E28').Value = '  -11.12%  '
$ws.Range('D29').Value = "'11.44"
$ws.Range('E29').Value = '  -10.26%  '
$ws.Range('D30').Value = "'8.50"
$ws.Range('E30').Value = '  -10.43%  '
$ws.Range('D31').Value = "'29.49"
$ws.Range('E31').Value = '  -7.37%  '
$ws.Range('D32').Value = "'662.56"
$ws.Range('E32').Value = '  -1.99%  '
$ws.Range('D33').Value = "'6.72"
$ws.Range('E33').Value = '  -16.47%  '
$ws.Range('B34').Value = 'Cosmos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D34').Value = "'11.17"
$ws.Range('E34').Value = '  -9.91%  '
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').Value = "'61.12"
$ws.Range('E35').Value = '  -4.37%  '
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('D38').Value = "'36.55"
$ws.Range('E38').Value = '  -13.72%  '
$ws.Range('D39').Value = "'0.380"
$ws.Range('E39').Value = '  -10.68%  '
$ws.Range('D40').Value = "'0.996"
$ws.Range('E40').Value = '  -0.23%  '
$ws.Range('D41').Value = "'0.127"
$ws.Range('E41').Value = '  -7.11%  '
$ws.Range('D42').Value = '2.831.16'
$ws.Range('E42').Value = '  -11.79%  '
$ws.Range('D43').Value = "'2.74"
$ws.Range('E43').Value = '  -16.67%  '
$ws.Range('D44').Value = '0.0₃0630'
$ws.Range('E44').Value = '  -22.36%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = "'0.0393"
$ws.Range('E45').Value = '  -7.08%  '
$ws.Range('B46').Value = 'WEMIXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').Value = "'2.61"
$ws.Range('E46').Value = '  -10.22%  '
$ws.Range('D47').Value = "'2.34"
$ws.Range('E47').Value = '  -14.63%  '
$ws.Range('D48').Value = "'137.48"
$ws.Range('E48').Value = '  -0.78%  '
$ws.Range('E49').Value = '  -6.37%  '
$ws.Range('D50').Value = "'2.82"
$ws.Range('E50').Value = '  -9.49%  '
$ws.Range('D51').Value = "'2.58"
$ws.Range('E51').Value = '  -6.53%  '
